$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (values use
# "." as a thousands separator, e.g. "66.910.10", and must not be
# auto-converted to numbers by Excel).

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '66.910.10'
$ws.Cells.Item(2, 5).Value = '  -2.37%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.477.53'
$ws.Cells.Item(3, 5).Value = '  -2.48%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '600.82'
$ws.Cells.Item(5, 5).Value = '  -3.23%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '147.49'
$ws.Cells.Item(6, 5).Value = '  -4.82%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.475.51'
$ws.Cells.Item(7, 5).Value = '  -2.40%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.479'
$ws.Cells.Item(9, 5).Value = '  -2.57%  '

$ws.Cells.Item(10, 5).Value = '  -3.30%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '7.61'
$ws.Cells.Item(11, 5).Value = '  +3.42%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.422'
$ws.Cells.Item(12, 5).Value = '  -3.69%  '

$ws.Cells.Item(13, 5).Value = '  -4.29%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.060.87'
$ws.Cells.Item(14, 5).Value = '  -2.62%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '31.31'
$ws.Cells.Item(15, 5).Value = '  -5.57%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.476.06'
$ws.Cells.Item(16, 5).Value = '  -2.55%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '66.911.65'
$ws.Cells.Item(17, 5).Value = '  -2.23%  '

$ws.Cells.Item(18, 5).Value = '  +0.20%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.41'
$ws.Cells.Item(19, 5).Value = '  -5.15%  '

$ws.Cells.Item(20, 5).Value = '  -4.60%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '10.02'
$ws.Cells.Item(21, 5).Value = '  +0.15%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '432.30'
$ws.Cells.Item(22, 5).Value = '  -5.01%  '

$ws.Cells.Item(23, 5).Value = '  -6.04%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '78.88'
$ws.Cells.Item(24, 5).Value = '  +0.24%  '

$ws.Cells.Item(25, 5).Value = '  +0.12%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '3.612.41'
$ws.Cells.Item(26, 5).Value = '  -2.62%  '

$ws.Cells.Item(27, 5).Value = '  -7.94%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '9.77'
$ws.Cells.Item(28, 5).Value = '  -7.54%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '8.35'
$ws.Cells.Item(29, 5).Value = '  -8.66%  '

$ws.Cells.Item(30, 5).Value = '  -3.48%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.58'
$ws.Cells.Item(31, 5).Value = '  -7.08%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.01'
$ws.Cells.Item(32, 5).Value = '  +0.89%  '

$ws.Cells.Item(33, 5).Value = '  -2.61%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '25.26'
$ws.Cells.Item(34, 5).Value = '  -3.50%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.464.78'
$ws.Cells.Item(35, 5).Value = '  -2.68%  '

$ws.Cells.Item(36, 5).Value = '  -6.84%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.79'
$ws.Cells.Item(37, 5).Value = '  -6.84%  '

$ws.Cells.Item(38, 5).Value = '  +0.01%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '7.87'
$ws.Cells.Item(39, 5).Value = '  -4.82%  '

$ws.Cells.Item(40, 5).Value = '  -0.18%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '173.95'
$ws.Cells.Item(41, 5).Value = '  -2.93%  '

$ws.Cells.Item(42, 5).Value = '  -4.47%  '

$ws.Cells.Item(43, 5).Value = '  -12.15%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '5.39'
$ws.Cells.Item(44, 5).Value = '  -3.78%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.894'
$ws.Cells.Item(45, 5).Value = '  -0.57%  '

$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '28.91'
$ws.Cells.Item(46, 5).Value = '  -6.69%  '

$ws.Cells.Item(47, 2).Value = 'OKB'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '46.32'
$ws.Cells.Item(47, 5).Value = '  +0.16%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.24'
$ws.Cells.Item(48, 5).Value = '  -7.90%  '

$ws.Cells.Item(49, 5).Value = '  -4.77%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.42'
$ws.Cells.Item(50, 5).Value = '  -9.14%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.975'
$ws.Cells.Item(51, 5).Value = '  -4.61%  '
